$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting so values like
# "29.103.74" or "0.00000000329" are not auto-coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.103.74'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").Value = '1.974.01'
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '329.52'
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("D6").Value = '1.009'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '0.4976'
$ws.Range("E7").Value = '  +0.96%  '
$ws.Range("E8").Value = '  +0.72%  '
$ws.Range("D9").Value = '53.17'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '0.09271'
$ws.Range("E10").Value = '  +4.95%  '
$ws.Range("D11").Value = '1.101'
$ws.Range("E11").Value = '  -1.18%  '
$ws.Range("D12").Value = '22.90'
$ws.Range("E12").Value = '  -1.13%  '
$ws.Range("D13").Value = '1.969.35'
$ws.Range("E13").Value = '  -3.75%  '
$ws.Range("D14").Value = '7.910'
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").Value = '6.463'
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("D16").Value = '1.011'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").Value = '0.00001109'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = '91.88'
$ws.Range("E18").Value = '  -4.41%  '
$ws.Range("D19").Value = '0.06726'
$ws.Range("E19").Value = '  +1.52%  '
$ws.Range("D20").Value = '19.28'
$ws.Range("E20").Value = '  -2.07%  '
$ws.Range("D21").Value = '1.008'
$ws.Range("D22").Value = '5.975'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").Value = '29.137.53'
$ws.Range("E23").Value = '  -1.17%  '
$ws.Range("D24").Value = '11.98'
$ws.Range("E24").Value = '  +1.01%  '
$ws.Range("D25").Value = '2.266'
$ws.Range("E25").Value = '  -0.94%  '
$ws.Range("D26").Value = '2.212.62'
$ws.Range("E26").Value = '  -3.06%  '
$ws.Range("E27").Value = '  +1.18%  '
$ws.Range("D28").Value = '155.31'
$ws.Range("E28").Value = '  -1.43%  '
$ws.Range("D29").Value = '6.367'
$ws.Range("E29").Value = '  -2.55%  '
$ws.Range("D30").Value = '2.268'
$ws.Range("E30").Value = '  -3.19%  '
$ws.Range("D31").Value = '126.97'
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("E32").Value = '  -0.20%  '
$ws.Range("D33").Value = '0.09854'
$ws.Range("D34").Value = '1.519'
$ws.Range("E34").Value = '  -2.35%  '
$ws.Range("D35").Value = '5.822'
$ws.Range("E35").Value = '  -0.23%  '
$ws.Range("D36").Value = '3.734'
$ws.Range("E36").Value = '  -0.99%  '
$ws.Range("D37").Value = '0.02429'
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("D38").Value = '1.323'
$ws.Range("E38").Value = '  +3.00%  '
$ws.Range("D39").Value = '0.06426'
$ws.Range("E39").Value = '  +0.75%  '
$ws.Range("D40").Value = '9.060'
$ws.Range("E40").Value = '  -5.49%  '
$ws.Range("D41").Value = '0.6486'
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").Value = '11.53'
$ws.Range("E42").Value = '  -1.86%  '
$ws.Range("D43").Value = '0.2003'
$ws.Range("E43").Value = '  -3.06%  '
$ws.Range("D44").Value = '1.008'
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '0.6223'
$ws.Range("E45").Value = '  -1.48%  '
$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '1.362'
$ws.Range("E46").Value = '  +7.25%  '
$ws.Range("D47").Value = '13.40'
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("E48").Value = '  -0.61%  '
$ws.Range("D49").Value = '3.482'
$ws.Range("E49").Value = '  -1.58%  '
$ws.Range("D50").Value = '0.00000000329'
$ws.Range("E50").Value = '  -2.17%  '
$ws.Range("D51").Value = '0.06969'
$ws.Range("E51").Value = '  -0.43%  '